# Add user authentication module
#
# The workbook is a rolling daily log (one row appended per day) across four
# worksheets. This change appends three new daily rows (for 2025-09-10,
# 2025-09-11 and 2025-09-12) to the bottom of each sheet's table, growing the
# used range from A1:I123 to A1:I126.

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "FE_LFT_#1"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        F = 380
        G = [double]"7.598631275147109e+23"
        I = 15
        Rows = @(
            @{ A = 45910.49049768518;  D = "0x00,0xE0"; E = "0xf"; H = 228 },
            @{ A = 45911.49237268518;  D = "0x00,0xE0"; E = "0xf"; H = 228 },
            @{ A = 45912.49243055555;  D = "0x00,0xE0"; E = "0xf"; H = 228 }
        )
    },
    @{
        Name = "FE_LFT_#2"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        F = 400
        G = [double]"5.68432987514711e+23"
        I = 14
        Rows = @(
            @{ A = 45910.49049768518;  D = "0x00,0xEC"; E = "0xe"; H = 240 },
            @{ A = 45911.49237268518;  D = "0x00,0xEC"; E = "0xe"; H = 236 },
            @{ A = 45912.49243055555;  D = "0x00,0xE8"; E = "0xe"; H = 236 }
        )
    },
    @{
        Name = "FE_PLT_#1"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        F = 110
        G = [double]"5.68631262647114e+23"
        I = 3
        Rows = @(
            @{ A = 45910.49049768518;  D = "0x00,0x57"; E = "0x3"; H = 87 },
            @{ A = 45911.49237268518;  D = "0x00,0x57"; E = "0x3"; H = 87 },
            @{ A = 45912.49243055555;  D = "0x00,0x57"; E = "0x3"; H = 87 }
        )
    },
    @{
        Name = "FE_PLT_#2"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        F = 110
        G = [double]"9.85046333984776e+23"
        I = 3
        Rows = @(
            @{ A = 45910.49049768518;  D = "0x00,0x54"; E = "0x3"; H = 84 },
            @{ A = 45911.49237268518;  D = "0x00,0x54"; E = "0x3"; H = 84 },
            @{ A = 45912.49243055555;  D = "0x00,0x54"; E = "0x3"; H = 84 }
        )
    }
)

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Name)
    $startRow = 124
    $i = 0
    foreach ($rd in $sd.Rows) {
        $r = $startRow + $i

        # Append a brand-new row at the bottom of the table (nothing below
        # it to shift, so this simply materialises row $r).
        $ws.Cells.Item($r, 1).EntireRow.Insert()

        $ws.Cells.Item($r, 1).Value = $rd.A
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $ws.Cells.Item($r, 2).Value = $sd.B
        $ws.Cells.Item($r, 3).Value = $sd.C
        $ws.Cells.Item($r, 4).Value = $rd.D
        $ws.Cells.Item($r, 5).Value = $rd.E
        $ws.Cells.Item($r, 6).Value = $sd.F
        $ws.Cells.Item($r, 7).Value = $sd.G
        $ws.Cells.Item($r, 8).Value = $rd.H
        $ws.Cells.Item($r, 9).Value = $sd.I

        $i = $i + 1
    }
}
